$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (no numeric auto-conversion) for cells whose literal representation
# (trailing zeros / percent signs) would otherwise be lost if Excel parsed them as numbers.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"

# Apply the updated values (diff: before.xlsx -> after.xlsx)
$ws.Range("D2").Value = "259.74"
$ws.Range("E2").Value = "5.98%"
$ws.Range("D3").Value = "27.92"
$ws.Range("E3").Value = "-2.41%"
$ws.Range("D4").Value = "5.212"
$ws.Range("E4").Value = "-0.59%"
$ws.Range("D5").Value = "0.05934"
$ws.Range("E5").Value = "4.18%"
$ws.Range("D6").Value = "6.736"
$ws.Range("E6").Value = "1.78%"
$ws.Range("D7").Value = "0.8733"
$ws.Range("E7").Value = "2.67%"
$ws.Range("D8").Value = "0.9972"
$ws.Range("E8").Value = "16.74%"
$ws.Range("E9").Value = "4.26%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.03636"
$ws.Range("E10").Value = "11.41%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07267"
$ws.Range("E11").Value = "2.48%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03230"
$ws.Range("E12").Value = "2.36%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09245"
$ws.Range("E13").Value = "0.47%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001540"
$ws.Range("E14").Value = "-0.07%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "0.0006039"
$ws.Range("E15").Value = "-93.98%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005782"
$ws.Range("E16").Value = "-3.53%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.496"
$ws.Range("E17").Value = "0.12%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "3.247"
$ws.Range("E18").Value = "1.88%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "2.210"
$ws.Range("E19").Value = "1.65%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3172"
$ws.Range("E20").Value = "0.16%"
$ws.Range("E21").Value = "-0.64%"
$ws.Range("D22").Value = "3.536"
$ws.Range("E22").Value = "1.50%"
$ws.Range("D23").Value = "0.04184"
$ws.Range("E23").Value = "2.59%"
$ws.Range("D24").Value = "0.1396"
$ws.Range("E24").Value = "1.31%"
$ws.Range("D25").Value = "0.001218"
$ws.Range("E25").Value = "-0.09%"
$ws.Range("D26").Value = "0.004575"
$ws.Range("E26").Value = "10.50%"
$ws.Range("D27").Value = "0.0001199"
$ws.Range("E27").Value = "-0.06%"
$ws.Range("D28").Value = "0.0001934"
$ws.Range("E28").Value = "33.55%"
$ws.Range("D40").Value = "0.03864"
$ws.Range("E40").Value = "2.94%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.005415"
$ws.Range("E41").Value = "4.48%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1112"
$ws.Range("E42").Value = "4.62%"
$ws.Range("D43").Value = "0.002380"
$ws.Range("E43").Value = "-4.39%"
$ws.Range("E44").Value = "16.47%"
$ws.Range("E45").Value = "2.75%"
$ws.Range("E46").Value = "-0.04%"
$ws.Range("D47").Value = "0.08534"
$ws.Range("E47").Value = "13.82%"
$ws.Range("D48").Value = "0.002138"
$ws.Range("E48").Value = "-12.32%"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").Value = "-0.04%"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").Value = "-0.04%"
